$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - VIN found -> green highlight + hyperlink to "10C029!B2"
$ws.Range("A2").Value = "KMHEC41BABA263951"
$ws.Range("B2").Value = "10C029!B2"

# Row 3 - VIN found -> green highlight + hyperlink to "10C029!B11"
$ws.Range("A3").Value = "KMHEC41CBBA240950"
$ws.Range("B3").Value = "10C029!B11"

# Row 4 - VIN not found -> red highlight, no hyperlink
$ws.Range("A4").Value = "KMH00000000000000"

# "Not found" style (red, indexed 52) applied first so it becomes fillId=3 / cellXfs index 1
$ws.Range("A4").Interior.PatternColorIndex = 52
$ws.Range("A4").Interior.Pattern = 9

# "Found" style (green, indexed 17) applied second so it becomes fillId=5 / cellXfs index 2
$ws.Range("A2").Interior.PatternColorIndex = 17
$ws.Range("A2").Interior.Pattern = 9
$ws.Range("A3").Interior.PatternColorIndex = 17
$ws.Range("A3").Interior.Pattern = 9

# Hyperlinks to the found VIN rows in sheet "10C029"
$ws.Hyperlinks.Add($ws.Range("B2"), "10C029!B2", "", "", "10C029")
$ws.Range("B2").Value = "10C029!B2"
$ws.Range("B2").ClearFormats()

$ws.Hyperlinks.Add($ws.Range("B3"), "10C029!B11", "", "", "10C029")
$ws.Range("B3").Value = "10C029!B11"
$ws.Range("B3").ClearFormats()

Write-Host "done"
